$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registered")
$ws.Activate()

# Replace the registered-user data row with the new user's details
$ws.Range("A2").Value = "738-94-4683"
$ws.Range("B2").Value = "Era"
$ws.Range("C2").Value = "Rice"
$ws.Range("D2").Value = "526 Gary Cape"
$ws.Range("E2").Value = "10539-6009"
$ws.Range("F2").Value = "Ernserhaven"
$ws.Range("G2").Value = "Georgia"
$ws.Range("H2").Value = "107-834-0930"
$ws.Range("I2").Value = "432-131-5405"
$ws.Range("J2").Value = "al.bradtke"
$ws.Range("K2").Value = "omer.gutmann@yahoo.com"
$ws.Range("L2").Value = "nQ6,LjR"

# Update the selection to span the full rows 2 through 4, active cell A2
$ws.Rows("2:4").Select()
